$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    40.100853000000001,
    39.946536999999999,
    40.53078,
    44.958874000000002,
    43.621464000000003,
    36.777602999999999,
    46.88138,
    37.578800000000001,
    43.728610000000003,
    39.494340000000001,
    54.792225000000002,
    43.195410000000003,
    35.122239999999998,
    47.962739999999997,
    41.124786,
    50.208480000000002,
    40.018830000000001,
    38.292313,
    49.403060000000004,
    43.034443000000003,
    44.354393000000002,
    36.866250000000001,
    42.497272000000002,
    43.42118,
    35.680250000000001,
    49.568171999999997,
    36.556429999999999,
    44.58229,
    45.114086,
    33.393639999999998,
    46.05442,
    35.413902,
    37.274470000000001,
    42.358919999999998,
    42.065849999999998,
    42.880443999999997,
    39.902287000000001,
    41.349105999999999,
    39.446002999999997,
    39.282955000000001,
    44.553513000000002,
    42.507862000000003,
    36.455886999999997,
    36.138522999999999,
    52.828957000000003,
    52.654130000000002,
    37.944476999999999,
    45.219917000000002,
    44.63306,
    45.073962999999999,
    39.933743,
    37.178420000000003,
    45.636054999999999,
    41.213290000000001,
    43.66507,
    47.574905000000001,
    31.628167999999999,
    50.437427999999997,
    42.888893000000003,
    42.397537,
    47.927599999999998,
    43.809401999999999,
    46.535310000000003,
    45.733800000000002,
    40.984406,
    38.535620000000002,
    39.898569999999999,
    46.461243000000003,
    45.819369999999999,
    36.224193999999997,
    45.800190000000001,
    41.541041999999997
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("F70").Select()

